# Add 4 new city rows to the KSA Cities worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Qassim region" label in the existing data uses a non-breaking space
# (U+00A0) between the two Arabic words, so reuse that exact string here to
# match the existing shared-string entry instead of creating a duplicate.
$qassimRegion = "منطقة" + [char]0x00A0 + "القصيم"

# Data to append: City, City English Name, City Arabic Name, Latitude, Longitude, Area, Region
$newRows = @(
    @("Al Bashayer",        "Al Bashayer",        "البشائر",        19.741637,  41.924493,  "منطقة عسير",      "جنوب المملكة"),
    @("Dhahran Al Janoub",  "Dhahran Al Janoub",  "ظهران الجنوب",   17.664974,  43.517555,  "منطقة عسير",      "جنوب المملكة"),
    @("Al Bathaa",          "Al Bathaa",          "البطحاء",        24.130382,  51.570082,  "المنطقة الشرقية", "شرق المملكة"),
    @("Oyoun Al Jawa",      "Oyoun Al Jawa",      "عيون الجواء",    26.516076,  43.619839,  $qassimRegion,     "وسط المملكة")
)

$startRow = 145
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]

    # Copy the style (border) from the row above so new rows look consistent.
    $ws.Range("A" + ($r - 1) + ":G" + ($r - 1)).Copy()
    $ws.Range("A" + $r + ":G" + $r).PasteSpecial(-4122) # xlPasteFormats
}

$ws.Range("A1:G" + ($startRow + $newRows.Count - 1)).Select()
